$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply scraped price (column D) and volume-change (column E) updates.
# Column D values must stay as literal text (matches source inlineStr cells),
# so we force a text number-format before assigning, then clear the format
# back off again so no stray style survives on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.161.59'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.668.56'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5118'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06405'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.59'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07415'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.671.69'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.514'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5815'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008581'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.22'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.224.83'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.938'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.82'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.48'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.205'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.32'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.624'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1191'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06361'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +13.00%  '
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.318'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.518'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.28%  '
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.014'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6075'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.371'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.59%  '
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.150'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.49%  '
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.082.78'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8645'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.10'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.817.61'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.20%  '
$ws.Range('E45').Value = '  +7.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.22'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.006'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.085'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05204'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4294'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.898'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.74%  '
